$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the formula in A1 to add 526
$ws.Range("A1").Formula = "=B1*C1+SUM(D1,E1,F1)/G1-H1+526"

# Reset selection back to A1 (default, removing explicit selection at A2)
$ws.Range("A1").Select()
